$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three "GasperiniShendure2019_*" rows (rows 11-13); remaining
# rows shift up to close the gap, matching the new dataset list.
$ws.Range("A11:A13").EntireRow.Delete() | Out-Null

# Add the new "Size (GB)" column (column E) with header formatting that
# matches the other header cells (bold, centered, bordered).
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 5).Value = "Size (GB)"

$sizeValues = @(
    0.03218394331634045,
    0.1295093791559339,
    0.4389201765879989,
    0.04276550840586424,
    0.467359134927392,
    0.03640089742839336,
    0.03129316493868828,
    0.2880413746461272,
    1.358732905238867,
    0.2302604019641876,
    1.359182251617312,
    0.6506966417655349,
    0.1371049117296934,
    0.0487448601052165,
    1.440504263155162,
    8.200729409232736,
    1.151940692216158,
    0.3528396841138601,
    1.528259848244488,
    0.2364100981503725,
    0.1352080181241035,
    0.2359374845400453,
    0.2504526963457465,
    0.3267531180754304,
    0.1441579991951585,
    0.269514954648912,
    0.2128301309421659,
    0.1092966264113784,
    0.5465821735560894
)

for ($i = 0; $i -lt $sizeValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $sizeValues[$i]
}
